$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;  E = 13.86384647080068;  G = 21.98653043760045 }
  3  = @{ B = 0.003078177322033415; C = 0.002658071450198252; D = 0.7210945179870265; E = 0.5333859586016987; G = 1.260216725360957 }
  4  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
  5  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
  6  = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
  7  = @{ B = 0.6545652718822623; C = 0.3048912486333797; D = 18.71679738969934; E = 0.5333859586016987; G = 20.20963986881668 }
  8  = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
  9  = @{ B = 0.6545652718822623; C = 0.3048912486333797; D = 0.7210945179870265; E = 13.86384647080068; G = 15.54439750930335 }
  10 = @{ B = 0.2881169905109251; C = 0.3048912486333797; D = 0.7210945179870265; E = 0.5333859586016987; G = 1.84748871573303 }
  11 = @{ B = 0.04172184405617529; C = 0.3048912486333797; D = 3.223369029078222; E = 13.86384647080068; G = 17.43382859256846 }
  12 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987; G = 5.582307763322248 }
  13 = @{ B = 0.6545652718822623; C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987; G = 3.536033448013082 }
  14 = @{ B = 0.2881169905109251; C = 0.3048912486333797; D = 3.223369029078222; E = 0.5333859586016987; G = 4.349763226824225 }
  15 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("B$row").Value = $vals.B
  $ws.Range("C$row").Value = $vals.C
  $ws.Range("D$row").Value = $vals.D
  $ws.Range("E$row").Value = $vals.E
  $ws.Range("G$row").Value = $vals.G
}
